$wb = $excel.ActiveWorkbook

# --- TestData sheet ---
# TC3/TC1 swap on rows 2 and 3, and Telemundo -> Oxygen on row 2
$wsTestData = $wb.Worksheets.Item("TestData")
$wsTestData.Range("A2").Value = "TC1"
$wsTestData.Range("A3").Value = "TC3"
$wsTestData.Range("E2").Value = "Oxygen"

# --- Windows sheet ---
# TC2 -> TC1 on rows 2 and 3, TC1 -> TC2 on row 4
$wsWindows = $wb.Worksheets.Item("Windows")
$wsWindows.Range("A2").Value = "TC1"
$wsWindows.Range("A3").Value = "TC1"
$wsWindows.Range("A4").Value = "TC2"

# Update the active selection shown on the Windows sheet (was D5, now A5)
[void]$wsWindows.Range("A5").Select()
